$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) need to be
# forced to Text format first, so Excel does not auto-convert them to
# numbers and silently drop significant trailing zeros (e.g. "0.06470").
$textCells = @("D4", "D5", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.270.07"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.676.65"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "217.22"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "0.2683"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "0.06470"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "0.07506"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "1.675.49"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "4.510"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "0.5765"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "0.000008497"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "64.58"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "26.313.66"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "4.906"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "189.75"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "6.180"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "144.91"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "7.801"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.1267"
$ws.Range("E26").Value = "  +6.35%  "
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "0.06492"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").Value = "1.362"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").Value = "1.317"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "3.583"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").Value = "3.587"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "1.653"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "0.6190"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").Value = "2.404"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "1.113.95"
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "0.8733"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").Value = "1.014"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "100.47"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "1.828.05"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "56.89"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "8.187"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05260"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.4291"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "6.060"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.3370"
$ws.Range("E51").Value = "  +1.94%  "
